$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.146235
$ws.Cells.Item(2, 8).Value = 0.438705
$ws.Cells.Item(2, 9).Value = 0.0224838618501081
$ws.Cells.Item(2, 10).Value = 0.0224838618501081
$ws.Cells.Item(2, 13).Value = 0.165747
$ws.Cells.Item(2, 14).Value = 0.497241
$ws.Cells.Item(2, 15).Value = 0.008095785894995438
$ws.Cells.Item(2, 16).Value = 0.00809578589499544
$ws.Cells.Item(2, 17).Value = 0.024238012545
$ws.Cells.Item(2, 18).Value = 0.218142112905
$ws.Cells.Item(2, 19).Value = 0.0001820245316311312
$ws.Cells.Item(2, 20).Value = 0.0001820245316311313

# Row 3
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.146235
$ws.Cells.Item(3, 8).Value = 0.438705
$ws.Cells.Item(3, 9).Value = 0.0224838618501081
$ws.Cells.Item(3, 10).Value = 0.0224838618501081
$ws.Cells.Item(3, 15).Value = 0.7079722685862583
$ws.Cells.Item(3, 16).Value = 0.7079722685862583
$ws.Cells.Item(3, 17).Value = 2.119601598915
$ws.Cells.Item(3, 18).Value = 19.076414390235
$ws.Cells.Item(3, 19).Value = 0.01591795068060106
$ws.Cells.Item(3, 20).Value = 0.01591795068060106

# Row 4
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.146235
$ws.Cells.Item(4, 8).Value = 0.438705
$ws.Cells.Item(4, 9).Value = 0.0224838618501081
$ws.Cells.Item(4, 10).Value = 0.0224838618501081
$ws.Cells.Item(4, 13).Value = 5.642879333333333
$ws.Cells.Item(4, 14).Value = 16.928638
$ws.Cells.Item(4, 15).Value = 0.2756221404547972
$ws.Cells.Item(4, 16).Value = 0.2756221404547972
$ws.Cells.Item(4, 17).Value = 0.82518645931
$ws.Cells.Item(4, 18).Value = 7.42667813379
$ws.Cells.Item(4, 19).Value = 0.00619705012881675
$ws.Cells.Item(4, 20).Value = 0.006197050128816751

# Row 5
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.146235
$ws.Cells.Item(5, 8).Value = 0.438705
$ws.Cells.Item(5, 9).Value = 0.0224838618501081
$ws.Cells.Item(5, 10).Value = 0.0224838618501081
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.1701286666666667
$ws.Cells.Item(5, 14).Value = 0.510386
$ws.Cells.Item(5, 15).Value = 0.008309805063949155
$ws.Cells.Item(5, 16).Value = 0.008309805063949155
$ws.Cells.Item(5, 17).Value = 0.02487876557
$ws.Cells.Item(5, 18).Value = 0.22390889013
$ws.Cells.Item(5, 19).Value = 0.0001868365090591615
$ws.Cells.Item(5, 20).Value = 0.0001868365090591615

# Row 6
$ws.Cells.Item(6, 7).Value = 3.793107666666666
$ws.Cells.Item(6, 9).Value = 0.5831962851568996
$ws.Cells.Item(6, 10).Value = 0.5831962851568997
$ws.Cells.Item(6, 13).Value = 0.165747
$ws.Cells.Item(6, 14).Value = 0.497241
$ws.Cells.Item(6, 15).Value = 0.008095785894995438
$ws.Cells.Item(6, 16).Value = 0.00809578589499544
$ws.Cells.Item(6, 17).Value = 0.6286962164269999
$ws.Cells.Item(6, 18).Value = 5.658265947843001
$ws.Cells.Item(6, 19).Value = 0.004721432259386965
$ws.Cells.Item(6, 20).Value = 0.004721432259386967

# Row 7
$ws.Cells.Item(7, 7).Value = 3.793107666666666
$ws.Cells.Item(7, 9).Value = 0.5831962851568996
$ws.Cells.Item(7, 10).Value = 0.5831962851568997
$ws.Cells.Item(7, 15).Value = 0.7079722685862583
$ws.Cells.Item(7, 16).Value = 0.7079722685862583
$ws.Cells.Item(7, 17).Value = 54.97915735031567
$ws.Cells.Item(7, 18).Value = 494.812416152841
$ws.Cells.Item(7, 19).Value = 0.4128867970336086
$ws.Cells.Item(7, 20).Value = 0.4128867970336086

# Row 8
$ws.Cells.Item(8, 7).Value = 3.793107666666666
$ws.Cells.Item(8, 9).Value = 0.5831962851568996
$ws.Cells.Item(8, 10).Value = 0.5831962851568997
$ws.Cells.Item(8, 13).Value = 5.642879333333333
$ws.Cells.Item(8, 14).Value = 16.928638
$ws.Cells.Item(8, 15).Value = 0.2756221404547972
$ws.Cells.Item(8, 16).Value = 0.2756221404547972
$ws.Cells.Item(8, 17).Value = 21.40404886134155
$ws.Cells.Item(8, 18).Value = 192.636439752074
$ws.Cells.Item(8, 19).Value = 0.1607418084202309
$ws.Cells.Item(8, 20).Value = 0.1607418084202309

# Row 9
$ws.Cells.Item(9, 7).Value = 3.793107666666666
$ws.Cells.Item(9, 9).Value = 0.5831962851568996
$ws.Cells.Item(9, 10).Value = 0.5831962851568997
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.1701286666666667
$ws.Cells.Item(9, 14).Value = 0.510386
$ws.Cells.Item(9, 15).Value = 0.008309805063949155
$ws.Cells.Item(9, 16).Value = 0.008309805063949155
$ws.Cells.Item(9, 17).Value = 0.6453163498531111
$ws.Cells.Item(9, 18).Value = 5.807847148677999
$ws.Cells.Item(9, 19).Value = 0.00484624744367314
$ws.Cells.Item(9, 20).Value = 0.00484624744367314

# Row 10
$ws.Cells.Item(10, 7).Value = 2.288493
$ws.Cells.Item(10, 8).Value = 6.865479000000001
$ws.Cells.Item(10, 9).Value = 0.3518594075080483
$ws.Cells.Item(10, 10).Value = 0.3518594075080483
$ws.Cells.Item(10, 13).Value = 0.165747
$ws.Cells.Item(10, 14).Value = 0.497241
$ws.Cells.Item(10, 15).Value = 0.008095785894995438
$ws.Cells.Item(10, 16).Value = 0.00809578589499544
$ws.Cells.Item(10, 17).Value = 0.3793108492710001
$ws.Cells.Item(10, 18).Value = 3.413797643439
$ws.Cells.Item(10, 19).Value = 0.002848578428325109
$ws.Cells.Item(10, 20).Value = 0.00284857842832511

# Row 11
$ws.Cells.Item(11, 7).Value = 2.288493
$ws.Cells.Item(11, 8).Value = 6.865479000000001
$ws.Cells.Item(11, 9).Value = 0.3518594075080483
$ws.Cells.Item(11, 10).Value = 0.3518594075080483
$ws.Cells.Item(11, 15).Value = 0.7079722685862583
$ws.Cells.Item(11, 16).Value = 0.7079722685862583
$ws.Cells.Item(11, 17).Value = 33.17053661507701
$ws.Cells.Item(11, 18).Value = 298.5348295356931
$ws.Cells.Item(11, 19).Value = 0.2491067029568897
$ws.Cells.Item(11, 20).Value = 0.2491067029568897

# Row 12
$ws.Cells.Item(12, 7).Value = 2.288493
$ws.Cells.Item(12, 8).Value = 6.865479000000001
$ws.Cells.Item(12, 9).Value = 0.3518594075080483
$ws.Cells.Item(12, 10).Value = 0.3518594075080483
$ws.Cells.Item(12, 13).Value = 5.642879333333333
$ws.Cells.Item(12, 14).Value = 16.928638
$ws.Cells.Item(12, 15).Value = 0.2756221404547972
$ws.Cells.Item(12, 16).Value = 0.2756221404547972
$ws.Cells.Item(12, 17).Value = 12.913689854178
$ws.Cells.Item(12, 18).Value = 116.223208687602
$ws.Cells.Item(12, 19).Value = 0.096980243036525
$ws.Cells.Item(12, 20).Value = 0.096980243036525

# Row 13
$ws.Cells.Item(13, 7).Value = 2.288493
$ws.Cells.Item(13, 8).Value = 6.865479000000001
$ws.Cells.Item(13, 9).Value = 0.3518594075080483
$ws.Cells.Item(13, 10).Value = 0.3518594075080483
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.1701286666666667
$ws.Cells.Item(13, 14).Value = 0.510386
$ws.Cells.Item(13, 15).Value = 0.008309805063949155
$ws.Cells.Item(13, 16).Value = 0.008309805063949155
$ws.Cells.Item(13, 17).Value = 0.3893382627660001
$ws.Cells.Item(13, 18).Value = 3.504044364894
$ws.Cells.Item(13, 19).Value = 0.002923883086308529
$ws.Cells.Item(13, 20).Value = 0.002923883086308529

# Row 14
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.2761626666666666
$ws.Cells.Item(14, 8).Value = 0.8284879999999999
$ws.Cells.Item(14, 9).Value = 0.04246044548494399
$ws.Cells.Item(14, 10).Value = 0.042460445484944
$ws.Cells.Item(14, 13).Value = 0.165747
$ws.Cells.Item(14, 14).Value = 0.497241
$ws.Cells.Item(14, 15).Value = 0.008095785894995438
$ws.Cells.Item(14, 16).Value = 0.00809578589499544
$ws.Cells.Item(14, 17).Value = 0.04577313351199999
$ws.Cells.Item(14, 18).Value = 0.411958201608
$ws.Cells.Item(14, 19).Value = 0.0003437506756522323
$ws.Cells.Item(14, 20).Value = 0.0003437506756522324

# Row 15
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.2761626666666666
$ws.Cells.Item(15, 8).Value = 0.8284879999999999
$ws.Cells.Item(15, 9).Value = 0.04246044548494399
$ws.Cells.Item(15, 10).Value = 0.042460445484944
$ws.Cells.Item(15, 15).Value = 0.7079722685862583
$ws.Cells.Item(15, 16).Value = 0.7079722685862583
$ws.Cells.Item(15, 17).Value = 4.002836734210666
$ws.Cells.Item(15, 18).Value = 36.025530607896
$ws.Cells.Item(15, 19).Value = 0.03006081791515895
$ws.Cells.Item(15, 20).Value = 0.03006081791515895

# Row 16
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.2761626666666666
$ws.Cells.Item(16, 8).Value = 0.8284879999999999
$ws.Cells.Item(16, 9).Value = 0.04246044548494399
$ws.Cells.Item(16, 10).Value = 0.042460445484944
$ws.Cells.Item(16, 13).Value = 5.642879333333333
$ws.Cells.Item(16, 14).Value = 16.928638
$ws.Cells.Item(16, 15).Value = 0.2756221404547972
$ws.Cells.Item(16, 16).Value = 0.2756221404547972
$ws.Cells.Item(16, 17).Value = 1.558352604371555
$ws.Cells.Item(16, 18).Value = 14.025173439344
$ws.Cells.Item(16, 19).Value = 0.01170303886922449
$ws.Cells.Item(16, 20).Value = 0.01170303886922449

# Row 17
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.2761626666666666
$ws.Cells.Item(17, 8).Value = 0.8284879999999999
$ws.Cells.Item(17, 9).Value = 0.04246044548494399
$ws.Cells.Item(17, 10).Value = 0.042460445484944
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.1701286666666667
$ws.Cells.Item(17, 14).Value = 0.510386
$ws.Cells.Item(17, 15).Value = 0.008309805063949155
$ws.Cells.Item(17, 16).Value = 0.008309805063949155
$ws.Cells.Item(17, 17).Value = 0.04698318626311111
$ws.Cells.Item(17, 18).Value = 0.4228486763679999
$ws.Cells.Item(17, 19).Value = 0.0003528380249083247
$ws.Cells.Item(17, 20).Value = 0.0003528380249083247
